$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation is inserted as row 147; all following rows
# (old 147..238) shift down to become rows 148..239.
$ws.Rows("147:147").Insert()

$ws.Range("A147").Value = 7
$ws.Range("B147").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C147").Value = "Ñuble"
$ws.Range("D147").Value = 44873
$ws.Range("E147").Value = 16
$ws.Range("F147").Value = 100112024
$ws.Range("G147").Value = "Choclo"
$ws.Range("H147").Value = "Dulce o Americano"
$ws.Range("I147").Value = "Primera"
$ws.Range("J147").Value = 50
$ws.Range("K147").Value = 30000
$ws.Range("L147").Value = 30000
$ws.Range("M147").Value = 30000
$ws.Range("N147").Value = "$/malla 70 unidades"
$ws.Range("O147").Value = "Región de Arica y Parinacota"
$ws.Range("P147").Value = 429
$ws.Range("Q147").Value = 70
$ws.Range("R147").Value = "Hortaliza"
